# BSS_RiskAndIssueList_V1.1.xlsx - "Update Risk and Issue"
#
# - Risk List: two risks that were "In Progress" are now "Closed".
# - Issue List: two issues that were "Active" are now "Fixed" (with a
#   Close_Date filled in), and one issue's Open_Date/Close_Date pair is
#   corrected.
# - The active worksheet moves from "Issue List" to "Risk List", and the
#   remembered selection on each of those two sheets is updated.

$wb = $excel.ActiveWorkbook

$riskList  = $wb.Worksheets.Item("Risk List")
$issueList = $wb.Worksheets.Item("Issue List")

# --- Risk List: mark two "In Progress" risks as "Closed" ---
$riskList.Range("E7").Value  = "Closed"
$riskList.Range("E18").Value = "Closed"

# --- Issue List: mark two "Active" issues as "Fixed" and set their Close_Date ---
$issueList.Range("E5").Value = "Fixed"
$issueList.Range("G5").Value = "1/3/2017"

$issueList.Range("E12").Value = "Fixed"
$issueList.Range("G12").Value = "1/3/2017"

# --- Issue List: correct the Open_Date / Close_Date on row 13 ---
$issueList.Range("F13").Value = "15/1/2017"
$issueList.Range("G13").Value = "1/2/2017"

# --- Update remembered selection on each sheet ---
$null = $riskList.Range("D10").Select()
$null = $issueList.Range("C12").Select()

# --- "Risk List" becomes the active sheet/tab (was "Issue List") ---
$riskList.Activate() | Out-Null
$null = $riskList.Range("D10").Select()
